# Tests 01 al 66 FacturaAfecta
# Appends 60 new "FacturaAfecta" test rows (FA_0001..FA_0049, FA_0056..FA_0066)
# to Hoja1, each with the fixed RUT "13712759-8" and software tag "Verity4.0",
# then freezes the header row / first column and scrolls/selects near the
# newly-added data (mirrors what a user does after pasting a big block of
# rows: freeze panes, and leave the selection on the last touched cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the ordered list of invoice codes: 0001-0049, then 0056-0066
# (the gap 0050-0055 is intentionally skipped - matches the source data).
$nums = @()
for ($i = 1; $i -le 49; $i++) { $nums += $i }
for ($i = 56; $i -le 66; $i++) { $nums += $i }

$startRow = 105
$row = $startRow
foreach ($n in $nums) {
    $code = "FA_{0:D4}" -f $n
    $ws.Cells.Item($row, 1).Value = $code
    $ws.Cells.Item($row, 2).Value = "13712759-8"
    $ws.Cells.Item($row, 3).Value = "Verity4.0"
    $row++
}
$lastRow = $row - 1

# Freeze the header row + first column, then land the selection/scroll
# near the bottom of the freshly pasted block.
$ws.Activate()
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

$excel.ActiveWindow.ScrollRow = $startRow + 45
$excel.ActiveWindow.ScrollColumn = 2

$ws.Range("F" + ($lastRow - 3)).Select()
